$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Remove the row for VarNum 2 (STARS_schnumb) entirely; all following rows shift up.
$ws.Rows.Item(3).Delete()

# Update the view's selection to match the saved state.
$ws.Range("E15").Select() | Out-Null
